$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D,E
$ws.Range("D2").Value = "34.139.61"
$ws.Range("E2").Value = "  +0.07%  "

# Row 3: update D,E
$ws.Range("D3").Value = "1.782.96"
$ws.Range("E3").Value = "  -0.36%  "

# Row 5: update D,E
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.09"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6: update D,E
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.547"
$ws.Range("E6").Value = "  +0.35%  "

# Row 7: update E
$ws.Range("E7").Value = "  +0.21%  "

# Row 8: update D,E
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.10"
$ws.Range("E8").Value = "  -0.66%  "

# Row 9: update E
$ws.Range("E9").Value = "  -0.82%  "

# Row 10: update D,E
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0686"
$ws.Range("E10").Value = "  -0.33%  "

# Row 11: update E
$ws.Range("E11").Value = "  +1.22%  "

# Row 12: update D,E
$ws.Range("D12").Value = "2.039.93"
$ws.Range("E12").Value = "  -0.38%  "

# Row 13: update B,C,D,E
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.99"
$ws.Range("E13").Value = "  -4.75%  "

# Row 14: update B,C,D,E
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.781.20"
$ws.Range("E14").Value = "  -0.64%  "

# Row 15: update D,E
$ws.Range("D15").Value = "34.131.76"
$ws.Range("E15").Value = "  +0.09%  "

# Row 16: update D,E
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.623"
$ws.Range("E16").Value = "  +0.12%  "

# Row 17: update E
$ws.Range("E17").Value = "  -0.02%  "

# Row 18: update D,E
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.62"
$ws.Range("E18").Value = "  -0.57%  "

# Row 19: update D,E
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.60"
$ws.Range("E19").Value = "  +0.67%  "

# Row 20: update D
$ws.Range("D20").Value = "0.0₃0788"

# Row 21: update B,C,D,E
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.93"
$ws.Range("E21").Value = "  -0.83%  "

# Row 22: update B,C,D,E
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.33%  "

# Row 23: update D,E
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.11"
$ws.Range("E23").Value = "  +0.20%  "

# Row 24: update E
$ws.Range("E24").Value = "  +0.10%  "

# Row 25: update D,E
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.16"
$ws.Range("E25").Value = "  +0.75%  "

# Row 26: update E
$ws.Range("E26").Value = "  -0.40%  "

# Row 27: update D,E
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.30"
$ws.Range("E27").Value = "  +0.26%  "

# Row 28: update E
$ws.Range("E28").Value = "  +0.62%  "

# Row 29: update D
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"

# Row 30: update B,C,D,E
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0519"
$ws.Range("E30").Value = "  +0.27%  "

# Row 31: update B,C,D,E
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.22"
$ws.Range("E31").Value = "  -0.78%  "

# Row 32: update D,E
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.71"
$ws.Range("E32").Value = "  +1.42%  "

# Row 33: update D,E
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.73"
$ws.Range("E33").Value = "  +3.13%  "

# Row 34: update D,E
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("E34").Value = "  -1.82%  "

# Row 35: update D,E
$ws.Range("D35").Value = "1.449.71"
$ws.Range("E35").Value = "  +3.29%  "

# Row 36: update D,E
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.50"
$ws.Range("E36").Value = "  +7.26%  "

# Row 37: update D,E
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.654"
$ws.Range("E37").Value = "  +0.48%  "

# Row 38: update E
$ws.Range("E38").Value = "  +1.28%  "

# Row 39: update E
$ws.Range("E39").Value = "  -0.30%  "

# Row 40: update D,E
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.41"
$ws.Range("E40").Value = "  +3.01%  "

# Row 41: update E
$ws.Range("E41").Value = "  +0.44%  "

# Row 42: update E
$ws.Range("E42").Value = "  -0.19%  "

# Row 43: update D,E
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.914"
$ws.Range("E43").Value = "  -0.81%  "

# Row 44: update D,E
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.67"
$ws.Range("E44").Value = "  +2.15%  "

# Row 45: update D,E
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0518"
$ws.Range("E45").Value = "  +2.43%  "

# Row 46: update B,C,D,E
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.06"
$ws.Range("E46").Value = "  +0.37%  "

# Row 47: update B,C,D,E
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.08"
$ws.Range("E47").Value = "  -0.06%  "

# Row 48: update D,E
$ws.Range("D48").Value = "1.939.57"
$ws.Range("E48").Value = "  -0.48%  "

# Row 49: update B,C,D,E
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.79"
$ws.Range("E49").Value = "  -2.03%  "

# Row 50: update B,C,D,E
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0131"
$ws.Range("E50").Value = "  -6.58%  "

# Row 51: update E
$ws.Range("E51").Value = "  +0.23%  "
